$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain decimal number need to be forced
# to text so Excel does not silently convert them to a Number cell (these
# columns store formatted price strings, not numeric values).
function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

$ws.Range('D2').Value = '34.348.36'
$ws.Range('E2').Value = '  +0.58%  '
$ws.Range('D3').Value = '1.835.64'
$ws.Range('E3').Value = '  +3.27%  '
Set-TextValue 'D4' '0.999'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('E5').Value = '  +0.27%  '
Set-TextValue 'D6' '0.558'
$ws.Range('E6').Value = '  +1.30%  '
Set-TextValue 'D7' '0.999'
$ws.Range('E7').Value = '  +0.00%  '
Set-TextValue 'D8' '32.06'
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('E9').Value = '  +4.33%  '
$ws.Range('E10').Value = '  +10.26%  '
Set-TextValue 'D11' '0.0932'
$ws.Range('E11').Value = '  +0.28%  '
$ws.Range('D12').Value = '2.101.01'
$ws.Range('E12').Value = '  +3.33%  '
$ws.Range('D13').Value = '1.832.75'
$ws.Range('E13').Value = '  +3.19%  '
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue 'D14' '10.83'
$ws.Range('E14').Value = '  -3.03%  '
$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue 'D15' '0.647'
$ws.Range('E15').Value = '  +3.25%  '
$ws.Range('D16').Value = '34.371.69'
$ws.Range('E16').Value = '  +0.74%  '
Set-TextValue 'D18' '69.96'
$ws.Range('E18').Value = '  +1.86%  '
Set-TextValue 'D19' '252.45'
$ws.Range('E19').Value = '  -0.96%  '
$ws.Range('D20').Value = '0.0₃0797'
$ws.Range('E20').Value = '  +7.89%  '
Set-TextValue 'D21' '11.22'
$ws.Range('E21').Value = '  +8.20%  '
$ws.Range('E22').Value = '  -0.17%  '
Set-TextValue 'D23' '4.29'
$ws.Range('E23').Value = '  +2.29%  '
$ws.Range('E24').Value = '  +1.23%  '
Set-TextValue 'D25' '160.48'
$ws.Range('E25').Value = '  +2.51%  '
$ws.Range('E26').Value = '  +2.23%  '
Set-TextValue 'D28' '0.115'
$ws.Range('E28').Value = '  +1.62%  '
$ws.Range('E29').Value = '  -0.09%  '
$ws.Range('E30').Value = '  +5.10%  '
Set-TextValue 'D31' '3.80'
$ws.Range('E31').Value = '  +0.84%  '
$ws.Range('E32').Value = '  +2.36%  '
$ws.Range('E33').Value = '  +0.15%  '
$ws.Range('E34').Value = '  +3.77%  '
$ws.Range('D35').Value = '1.450.92'
$ws.Range('E35').Value = '  +0.89%  '
Set-TextValue 'D36' '0.649'
$ws.Range('E36').Value = '  +4.10%  '
Set-TextValue 'D37' '1.07'
$ws.Range('E37').Value = '  +1.84%  '
$ws.Range('E38').Value = '  +3.05%  '
$ws.Range('E39').Value = '  +9.45%  '
Set-TextValue 'D40' '82.17'
$ws.Range('E40').Value = '  -0.70%  '
$ws.Range('E41').Value = '  -2.86%  '
$ws.Range('E42').Value = '  +0.39%  '
Set-TextValue 'D43' '2.15'
$ws.Range('E43').Value = '  +4.63%  '
Set-TextValue 'D44' '6.11'
$ws.Range('E44').Value = '  +4.97%  '
$ws.Range('D45').Value = '1.997.19'
$ws.Range('E45').Value = '  +3.17%  '
Set-TextValue 'D46' '0.0500'
$ws.Range('E46').Value = '  -2.18%  '
$ws.Range('E47').Value = '  +0.30%  '
$ws.Range('E48').Value = '  +8.88%  '
$ws.Range('E49').Value = '  -0.04%  '
Set-TextValue 'D50' '11.93'
$ws.Range('E50').Value = '  -2.12%  '
$ws.Range('E51').Value = '  +5.98%  '
